# Generate Report for Handoff
#
# A new source file (ffffaf8b89a2-1303-4afe-8171-aa7f6d5d2396.md) has been
# added to the localization set, and the existing source file's generated
# UUID changed from 94f631a3-6c80-40bb-b9f5-46cbfbfaa00e to
# 42de5bca-bee6-4481-acdf-58d06ea59353 (new handoff package hash
# 2876460dd14dd29860c1a9c1343044e64ae3d965) with updated handoff timestamps.
# This pushes the ".localization-config" row down by one on every sheet.

$wb = $excel.ActiveWorkbook

$oldUuid = "94f631a3-6c80-40bb-b9f5-46cbfbfaa00e"
$newUuid = "42de5bca-bee6-4481-acdf-58d06ea59353"
$newFileUuid = "ffffaf8b89a2-1303-4afe-8171-aa7f6d5d2396"
$pkgHash = "2876460dd14dd29860c1a9c1343044e64ae3d965"

$newMdName = "$newUuid.md"
$newFileMdName = "$newFileUuid.md"
$zhXlfName = "$newUuid.$pkgHash.zh-cn.xlf"
$deXlfName = "$newUuid.$pkgHash.de-de.xlf"

$zhHandoffTime = "2016-03-02 15:25:09"
$deHandoffTime = "2016-03-02 15:25:20"
$epoch = "0001-01-01 00:00:00"

$newMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/0b733104bb2dc6b49bd24adad44fe11b9a6642ae/e2e/$newMdName"
$newFileMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/0b733104bb2dc6b49bd24adad44fe11b9a6642ae/e2e/$newFileMdName"
$localizationConfigUrl = "https://github.com/OpenLocalizationTest/oltest/blob/0b733104bb2dc6b49bd24adad44fe11b9a6642ae/.localization-config"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1fc032d51dd6fabcc87187101b80b4b1611449cb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlfName"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/85ff2edfb74776d20974146677c01df68aad7bce/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlfName"

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Drop the existing hyperlinks (on A2/A3) up front -- every hyperlinked
# cell is being rewritten below, and this engine's Hyperlinks.Delete()
# only operates cleanly at the whole-collection level.
$ov.Hyperlinks.Delete()

$ov.Range("A2").Value = $newMdName
$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"

$ov.Range("A3").Value = $newFileMdName
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"

$ov.Range("A4").Value = ".localization-config"
$ov.Range("B4").Value = "Not to be localized"
$ov.Range("C4").Value = "Not to be localized"

$ov.Range("A2").Style = "HyperLink"
$ov.Range("A3").Style = "HyperLink"
$ov.Range("A4").Style = "HyperLink"

$ov.Hyperlinks.Add($ov.Range("A2"), $newMdUrl, [Type]::Missing, [Type]::Missing, $newMdName)
$ov.Hyperlinks.Add($ov.Range("A3"), $newFileMdUrl, [Type]::Missing, [Type]::Missing, $newFileMdName)
$ov.Hyperlinks.Add($ov.Range("A4"), $localizationConfigUrl, [Type]::Missing, [Type]::Missing, ".localization-config")

# ---------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | Status | Latest Handoff File |
#   Latest Handoff Datetime | Latest Target File | Latest Handback File |
#   Latest Handback DateTime | Handoff Reason | Dependency From
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Hyperlinks.Delete()

$zh.Range("A2").Value = $newMdName
$zh.Range("B2").Value = "Ready for handoff"
$zh.Range("C2").Value = $zhXlfName
$zh.Range("D2").Value = $zhHandoffTime
$zh.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("G2").Value = $epoch
$zh.Range("H2").Value = "Include"

$zh.Range("A3").Value = $newFileMdName
$zh.Range("B3").Value = "Ready for handoff"
$zh.Range("C3").Value = $zhXlfName
$zh.Range("D3").Value = $zhHandoffTime
$zh.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("G3").Value = $epoch
$zh.Range("H3").Value = "Include"

$zh.Range("A4").Value = ".localization-config"
$zh.Range("B4").Value = "Not to be localized"
$zh.Range("D4").Value = $epoch
$zh.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("G4").Value = $epoch
$zh.Range("H4").Value = "Ignored"

$zh.Range("A2").Style = "HyperLink"
$zh.Range("C2").Style = "HyperLink"
$zh.Range("A3").Style = "HyperLink"
$zh.Range("C3").Style = "HyperLink"
$zh.Range("A4").Style = "HyperLink"

$zh.Hyperlinks.Add($zh.Range("A2"), $newMdUrl, [Type]::Missing, [Type]::Missing, $newMdName)
$zh.Hyperlinks.Add($zh.Range("C2"), $zhXlfUrl, [Type]::Missing, [Type]::Missing, $zhXlfName)
$zh.Hyperlinks.Add($zh.Range("A3"), $newFileMdUrl, [Type]::Missing, [Type]::Missing, $newFileMdName)
$zh.Hyperlinks.Add($zh.Range("C3"), $zhXlfUrl, [Type]::Missing, [Type]::Missing, $zhXlfName)
$zh.Hyperlinks.Add($zh.Range("A4"), $localizationConfigUrl, [Type]::Missing, [Type]::Missing, ".localization-config")

# ---------------------------------------------------------------------
# Sheet "de-de": same layout as "zh-cn"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Hyperlinks.Delete()

$de.Range("A2").Value = $newMdName
$de.Range("B2").Value = "Ready for handoff"
$de.Range("C2").Value = $deXlfName
$de.Range("D2").Value = $deHandoffTime
$de.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("G2").Value = $epoch
$de.Range("H2").Value = "Include"

$de.Range("A3").Value = $newFileMdName
$de.Range("B3").Value = "Ready for handoff"
$de.Range("C3").Value = $deXlfName
$de.Range("D3").Value = $deHandoffTime
$de.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("G3").Value = $epoch
$de.Range("H3").Value = "Include"

$de.Range("A4").Value = ".localization-config"
$de.Range("B4").Value = "Not to be localized"
$de.Range("D4").Value = $epoch
$de.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("G4").Value = $epoch
$de.Range("H4").Value = "Ignored"

$de.Range("A2").Style = "HyperLink"
$de.Range("C2").Style = "HyperLink"
$de.Range("A3").Style = "HyperLink"
$de.Range("C3").Style = "HyperLink"
$de.Range("A4").Style = "HyperLink"

$de.Hyperlinks.Add($de.Range("A2"), $newMdUrl, [Type]::Missing, [Type]::Missing, $newMdName)
$de.Hyperlinks.Add($de.Range("C2"), $deXlfUrl, [Type]::Missing, [Type]::Missing, $deXlfName)
$de.Hyperlinks.Add($de.Range("A3"), $newFileMdUrl, [Type]::Missing, [Type]::Missing, $newFileMdName)
$de.Hyperlinks.Add($de.Range("C3"), $deXlfUrl, [Type]::Missing, [Type]::Missing, $deXlfName)
$de.Hyperlinks.Add($de.Range("A4"), $localizationConfigUrl, [Type]::Missing, [Type]::Missing, ".localization-config")

"Handoff report regenerated"
